$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.774.48'
$ws.Range('E2').Value = '  -0.41%  '

# Row 3
$ws.Range('D3').Value = '1.635.45'
$ws.Range('E3').Value = '  -0.30%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.22%  '

# Row 6
$ws.Range('E6').Value = '  -0.35%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('E8').Value = '  -0.22%  '

# Row 9
$ws.Range('E9').Value = '  -0.54%  '

# Row 10
$ws.Range('E10').Value = '  +0.66%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.45%  '

# Row 12
$ws.Range('E12').Value = '  -1.01%  '

# Row 13
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.860.69'
$ws.Range('E13').Value = '  -0.34%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.628.04'
$ws.Range('E14').Value = '  -0.85%  '

# Row 15
$ws.Range('E15').Value = '  -1.20%  '

# Row 16
$ws.Range('D16').Value = '0.0₃0775'
$ws.Range('E16').Value = '  +1.60%  '

# Row 17
$ws.Range('E17').Value = '  -0.06%  '

# Row 18
$ws.Range('D18').Value = '25.793.89'
$ws.Range('E18').Value = '  -0.40%  '

# Row 19
$ws.Range('E19').Value = '  -0.04%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.65%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.89%  '

# Row 22
$ws.Range('E22').Value = '  +0.34%  '

# Row 23
$ws.Range('E23').Value = '  +0.06%  '

# Row 24
$ws.Range('E24').Value = '  +0.05%  '

# Row 25
$ws.Range('E25').Value = '  -0.99%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.35'
$ws.Range('D26').Style = 'Normal'

# Row 27
$ws.Range('E27').Value = '  -5.43%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.30%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.42%  '

# Row 30
$ws.Range('E30').Value = '  -0.06%  '

# Row 31
$ws.Range('E31').Value = '  +0.41%  '

# Row 32
$ws.Range('E32').Value = '  +1.36%  '

# Row 33
$ws.Range('E33').Value = '  +0.99%  '

# Row 34
$ws.Range('E34').Value = '  +1.31%  '

# Row 35
$ws.Range('E35').Value = '  +0.45%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.897'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.78%  '

# Row 37
$ws.Range('E37').Value = '  -0.20%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.551'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.65%  '

# Row 39
$ws.Range('D39').Value = '1.109.92'
$ws.Range('E39').Value = '  -1.68%  '

# Row 40
$ws.Range('E40').Value = '  +0.21%  '

# Row 41
$ws.Range('E41').Value = '  +0.00%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.27%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.804'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.10%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.19%  '

# Row 45
$ws.Range('D45').Value = '0.0₆0113'
$ws.Range('E45').Value = '  +2.45%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.33%  '

# Row 47
$ws.Range('E47').Value = '  +11.68%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.72'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.28%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0503'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.44%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
